$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Tumakuru (Tumkur)"
$ws.Range("G9").Value = "Ballari (Bellary)"
$ws.Range("G10").Value = "Ballari (Bellary)"
$ws.Range("G14").Value = "Ballari (Bellary)"
$ws.Range("G21").Value = "Uttara Kannada (Karwar)"
$ws.Range("G31").Value = "Tumakuru (Tumkur)"
$ws.Range("G35").Value = "Ballari (Bellary)"
$ws.Range("G36").Value = "Vijayapura (Bijapur)"
$ws.Range("G37").Value = "Ballari (Bellary)"
$ws.Range("G41").Value = "Ballari (Bellary)"
